$wb = $excel.ActiveWorkbook

# Sheet: ALC (index 1)
$ws = $wb.Worksheets.Item(1)
# Row 17
$ws.Range("H17").Value = 4623.2856
$ws.Range("J17").Value = 3966.6667
$ws.Range("L17").Value = 11900.0001
$ws.Range("N17").Value = -12236.0001

# Row 28
$ws.Range("H28").Value = 194.14285
$ws.Range("I28").Value = 194.14285
$ws.Range("K28").Value = 194.14285
$ws.Range("M28").Value = 290.85715

# Row 62
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

# Row 65
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

# Row 86
$ws.Range("H86").Value = 1999
$ws.Range("I86").Value = 1999
$ws.Range("K86").Value = 1999
$ws.Range("M86").Value = -876

# Row 89
$ws.Range("H89").Value = 1999
$ws.Range("I89").Value = 1999
$ws.Range("K89").Value = 9995
$ws.Range("M89").Value = -4379

# Row 112
$ws.Range("H112").Value = 2965.2144
$ws.Range("J112").Value = 3501.182
$ws.Range("L112").Value = 10503.546
$ws.Range("N112").Value = -12719.546

# Row 121
$ws.Range("H121").Value = 854
$ws.Range("J121").Value = 877.3570999999999
$ws.Range("L121").Value = 2632.0713
$ws.Range("N121").Value = -6126.0713

# Row 125
$ws.Range("H125").Value = 520.3
$ws.Range("J125").Value = 478.33334
$ws.Range("L125").Value = 4305.00006
$ws.Range("N125").Value = -9225.00006

# Row 137
$ws.Range("H137").Value = 37072.43
$ws.Range("I137").Value = 1431
$ws.Range("J137").Value = 67961.664
$ws.Range("K137").Value = 4293
$ws.Range("L137").Value = 203884.992
$ws.Range("M137").Value = -1743
$ws.Range("N137").Value = -208984.992

# Row 138
$ws.Range("H138").Value = 4892.648
$ws.Range("J138").Value = 4650.912
$ws.Range("L138").Value = 13952.736
$ws.Range("N138").Value = -24232.736

# Sheet: ARM (index 2)
$ws = $wb.Worksheets.Item(2)
# Row 4
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -282

# Row 61
$ws.Range("H61").Value = 20397.361
$ws.Range("I61").Value = 26743.469
$ws.Range("J61").Value = 6859
$ws.Range("K61").Value = 26743.469
$ws.Range("L61").Value = 6859
$ws.Range("M61").Value = -26531.469
$ws.Range("N61").Value = -7283

# Row 74
$ws.Range("H74").Value = 2150
$ws.Range("I74").Value = 733.3333
$ws.Range("K74").Value = 733.3333
$ws.Range("M74").Value = 140.6667

# Row 77
$ws.Range("H77").Value = 2150
$ws.Range("I77").Value = 733.3333
$ws.Range("K77").Value = 3666.6665
$ws.Range("M77").Value = 701.3334999999997

# Row 97
$ws.Range("H97").Value = 2257.4285
$ws.Range("I97").Value = 1701.3077
$ws.Range("J97").Value = 3161.125
$ws.Range("K97").Value = 1701.3077
$ws.Range("L97").Value = 3161.125
$ws.Range("M97").Value = -1205.3077
$ws.Range("N97").Value = -4153.125

# Row 132
$ws.Range("H132").Value = 2039.1346
$ws.Range("I132").Value = 1797.0385
$ws.Range("K132").Value = 5391.1155
$ws.Range("M132").Value = -2861.1155

# Row 136
$ws.Range("H136").Value = 20397.361
$ws.Range("I136").Value = 26743.469
$ws.Range("J136").Value = 6859
$ws.Range("K136").Value = 80230.40700000001
$ws.Range("L136").Value = 20577
$ws.Range("M136").Value = -77680.40700000001
$ws.Range("N136").Value = -25677

# Sheet: BSM (index 3)
$ws = $wb.Worksheets.Item(3)
# Row 134
$ws.Range("H134").Value = 8817.096
$ws.Range("I134").Value = 9281.056
$ws.Range("K134").Value = 27843.168
$ws.Range("M134").Value = -25308.168

# Sheet: CRP (index 4)
$ws = $wb.Worksheets.Item(4)
# Row 58
$ws.Range("H58").Value = 2175888.5
$ws.Range("I58").Value = 2900235.5
$ws.Range("K58").Value = 2900235.5
$ws.Range("M58").Value = -2900032.5

# Row 100
$ws.Range("H100").Value = 40000
$ws.Range("J100").Value = 40000
$ws.Range("L100").Value = 40000
$ws.Range("N100").Value = -42164

# Row 132
$ws.Range("H132").Value = 2243.5
$ws.Range("I132").Value = 1617.4166
$ws.Range("K132").Value = 4852.2498
$ws.Range("M132").Value = -2322.2498

# Row 134
$ws.Range("H134").Value = 1153.9048
$ws.Range("I134").Value = 1026.5938
$ws.Range("K134").Value = 3079.7814
$ws.Range("M134").Value = -544.7814000000003

# Row 136
$ws.Range("H136").Value = 2175888.5
$ws.Range("I136").Value = 2900235.5
$ws.Range("K136").Value = 8700706.5
$ws.Range("M136").Value = -8698156.5

# Sheet: CUL (index 5)
$ws = $wb.Worksheets.Item(5)
# Row 5
$ws.Range("H5").Value = 711.0476
$ws.Range("J5").Value = 918
$ws.Range("L5").Value = 2754
$ws.Range("N5").Value = -2978

# Row 68
$ws.Range("H68").Value = 500000
$ws.Range("J68").Value = 500000
$ws.Range("L68").Value = 1500000
$ws.Range("N68").Value = -1501622

# Row 71
$ws.Range("H71").Value = 500000
$ws.Range("J71").Value = 500000
$ws.Range("L71").Value = 4500000
$ws.Range("N71").Value = -4508112

# Row 86
$ws.Range("H86").Value = 500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = -3872

# Row 89
$ws.Range("H89").Value = 500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = -16356

# Row 130
$ws.Range("H130").Value = 2353
$ws.Range("I130").Value = 1423.4
$ws.Range("K130").Value = 4270.200000000001
$ws.Range("M130").Value = 749.7999999999993

# Row 131
$ws.Range("H131").Value = 11541.703
$ws.Range("J131").Value = 12006
$ws.Range("L131").Value = 36018
$ws.Range("N131").Value = -46098

# Row 132
$ws.Range("H132").Value = 1249.8
$ws.Range("I132").Value = 949.6667
$ws.Range("K132").Value = 8547.0003
$ws.Range("M132").Value = -6017.0003

# Row 134
$ws.Range("H134").Value = 2361
$ws.Range("I134").Value = 2089.6086
$ws.Range("J134").Value = 2985.2
$ws.Range("K134").Value = 6268.825800000001
$ws.Range("L134").Value = 8955.599999999999
$ws.Range("M134").Value = -1198.825800000001
$ws.Range("N134").Value = -19095.6

# Row 135
$ws.Range("H135").Value = 711.0476
$ws.Range("J135").Value = 918
$ws.Range("L135").Value = 8262
$ws.Range("N135").Value = -13332

# Sheet: GSM (index 6)
$ws = $wb.Worksheets.Item(6)
# Row 80
$ws.Range("H80").Value = 2999
$ws.Range("I80").Value = 2999
$ws.Range("K80").Value = 2999
$ws.Range("M80").Value = -2001

# Row 83
$ws.Range("H83").Value = 2999
$ws.Range("I83").Value = 2999
$ws.Range("K83").Value = 14995
$ws.Range("M83").Value = -10003

# Row 95
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -20492

# Row 122
$ws.Range("H122").Value = 1976
$ws.Range("I122").Value = 1952.5
$ws.Range("K122").Value = 5857.5
$ws.Range("M122").Value = -3407.5

# Row 126
$ws.Range("H126").Value = 1716626
$ws.Range("I126").Value = 2225069.8
$ws.Range("K126").Value = 6675209.399999999
$ws.Range("M126").Value = -6672739.399999999

# Row 132
$ws.Range("H132").Value = 1546045.5
$ws.Range("I132").Value = 2413148
$ws.Range("K132").Value = 7239444
$ws.Range("M132").Value = -7236914

# Sheet: LTW (index 7)
$ws = $wb.Worksheets.Item(7)
# Row 68
$ws.Range("H68").Value = 2333
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 2999
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 2999
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -4497

# Row 71
$ws.Range("H71").Value = 2333
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 2999
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 14995
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -22483

# Row 136
$ws.Range("H136").Value = 4527.4443
$ws.Range("I136").Value = 3549.4
$ws.Range("K136").Value = 10648.2
$ws.Range("M136").Value = -8098.200000000001

# Sheet: WVR (index 8)
$ws = $wb.Worksheets.Item(8)
# Row 132
$ws.Range("H132").Value = 1389.1052
$ws.Range("I132").Value = 1024.3125
$ws.Range("K132").Value = 3072.9375
$ws.Range("M132").Value = -542.9375

# Row 136
$ws.Range("H136").Value = 42739384
$ws.Range("I136").Value = 79368140
$ws.Range("J136").Value = 5830.8335
$ws.Range("K136").Value = 238104420
$ws.Range("L136").Value = 17492.5005
$ws.Range("M136").Value = -238101870
$ws.Range("N136").Value = -22592.5005
